$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Cells.Item(1,1)
$x = $c.Font.Charset
Write-Host "VALUE:[$x]"
Write-Host $x.GetType()
